# Roadmap.xlsx update
# [Engine] [Render] refactor cameraController
#
# Sheet "路线图总览" (Worksheets.Item(1)):
#   - selection moves from B4 to B15 (no data change)
# Sheet "Roadmap" (Worksheets.Item(2)):
#   - row 8 (Forward Pass blin-phong row) gets marked "in progress" (green fill) and
#     gains a start date 2023.11.12 + progress marker "进行中"
#   - row 11 (调整RenderCameraController) switches from the "in progress" (green)
#     fill to the "done" (blue) fill and gains an end date 2023.11.12
#   - row 12 is a brand new task: Renderer / 调整Material
#   - row 13 is a brand new task: Editor / 整理资源目录结构
#   - selection moves from E13 to C11

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# A scratch cell used to coerce date-shaped strings ("2023.11.12", ...) into plain
# text instead of having them auto-parsed into date serials: set the scratch cell
# to text format, give it the string we want, then copy only the *value* over to
# the real target (which already carries the fill/style we want).
function Set-TextValue($range, $text) {
    $scratch = $ws2.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163) # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# Sheet2 "Roadmap" content changes
# ---------------------------------------------------------------------------

# Row 8: apply the "in progress" (green) style used by row 7 to the whole row,
# then fill in the newly tracked start date + progress marker.
$ws2.Range("A7:C7").Copy()
$ws2.Range("A8:C8").PasteSpecial(-4122) # xlPasteFormats
$ws2.Range("G7").Copy()
$ws2.Range("G8").PasteSpecial(-4122)
$ws2.Range("E7").Copy()
$ws2.Range("E8").PasteSpecial(-4122)

Set-TextValue $ws2.Range("E8") "2023.11.12"
$ws2.Range("G8").Value = "进行中"

# Row 11: "调整RenderCameraController" switches to the "done" (blue) fill used by
# rows 9/10, and records a new end date in F11.
$ws2.Range("A9:C9").Copy()
$ws2.Range("A11:C11").PasteSpecial(-4122)
$ws2.Range("E9").Copy()
$ws2.Range("E11").PasteSpecial(-4122)
$ws2.Range("G9").Copy()
$ws2.Range("G11").PasteSpecial(-4122)
$ws2.Range("F9").Copy()
$ws2.Range("F11").PasteSpecial(-4122)

Set-TextValue $ws2.Range("F11") "2023.11.12"

# Row 12 (new task): Renderer / 调整Material, still "in progress" (green fill).
$ws2.Range("A7:D7").Copy()
$ws2.Range("A12:D12").PasteSpecial(-4122)
$ws2.Range("E7").Copy()
$ws2.Range("E12").PasteSpecial(-4122)
$ws2.Range("G7").Copy()
$ws2.Range("G12").PasteSpecial(-4122)

$ws2.Range("A12").Value = "Sean Duan"
$ws2.Range("B12").Value = "Renderer"
$ws2.Range("C12").Value = "调整Material"
$ws2.Range("D12").Value = "根据重构的管线，修改Material反射"
Set-TextValue $ws2.Range("E12") "2023.11.12"
$ws2.Range("G12").Value = "进行中"

# Row 13 (new task): Editor / 整理资源目录结构, plain/unstyled (same as rows 1/6).
$ws2.Range("A13").Value = "Sean Duan"
$ws2.Range("B13").Value = "Editor"
$ws2.Range("C13").Value = "整理资源目录结构"

# ---------------------------------------------------------------------------
# Selections / active cells
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B15").Select()

$ws2.Activate()
$ws2.Range("C11").Select()
